$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.584.99"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "'1.912.89"
$ws.Range("E3").Value = "  +5.19%  "
$ws.Range("D5").Value = "'314.77"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.5193"
$ws.Range("E7").Value = "  +4.50%  "
$ws.Range("D8").Value = "'0.3968"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").Value = "'0.09764"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E10").Value = "  +4.12%  "
$ws.Range("D11").Value = "'41.98"
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("D12").Value = "'6.557"
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("D13").Value = "'21.25"
$ws.Range("E13").Value = "  +3.37%  "
$ws.Range("D14").Value = "'1.908.66"
$ws.Range("E14").Value = "  +5.20%  "
$ws.Range("D15").Value = "'7.538"
$ws.Range("E15").Value = "  +3.54%  "
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'94.81"
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'0.06645"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'18.21"
$ws.Range("E20").Value = "  +5.99%  "
$ws.Range("D21").Value = "'0.9999"
$ws.Range("D22").Value = "'6.328"
$ws.Range("E22").Value = "  +5.85%  "
$ws.Range("D23").Value = "'28.688.90"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").Value = "'11.58"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").Value = "'2.300"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.686"
$ws.Range("E26").Value = "  +11.45%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "'2.130.19"
$ws.Range("E27").Value = "  +5.23%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'21.32"
$ws.Range("E28").Value = "  +2.68%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'158.59"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'128.91"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.117"
$ws.Range("E31").Value = "  +7.78%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1085"
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.774"
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.634"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "'9.976"
$ws.Range("E35").Value = "  +12.33%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06797"
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02435"
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.264"
$ws.Range("E38").Value = "  +7.29%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2234"
$ws.Range("E39").Value = "  +4.34%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'11.84"
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6496"
$ws.Range("E41").Value = "  +4.63%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "'5.086"
$ws.Range("E42").Value = "  +2.47%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.191"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.59"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6116"
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.763"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.282"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.041"
$ws.Range("E49").Value = "  +5.20%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'125.17"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.209"
$ws.Range("E51").Value = "  +2.49%  "
